$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu")

$ws.Cells.Item(6, 1).Value = "Fried RiceFried Rice"
$ws.Cells.Item(6, 2).Value = 2400

$ws.Cells.Item(7, 1).Value = "Fried Rice"
$ws.Cells.Item(7, 2).Value = 1200
